$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 440, shifting the existing rows 440:476 down to 443:479.
$xlShiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$ws.Rows("440:442").Insert($xlShiftDown)

# New rows (440:442) share the same constant columns as the rest of the
# "Betarraga" block: A, B, C, E, F, G, H, N, O, Q, R.
$constA = 6
$constB = "Mercado Mayorista Lo Valledor de Santiago"
$constC = "Metropolitana"
$constE = 13
$constF = 100114014
$constG = "Betarraga"
$constH = "Sin especificar"
$constN = "`$/unidad"
$constO = "Región Metropolitana"
$constQ = 1
$constR = "Hortaliza"

$newRows = @(
    @{ Row = 440; Fecha = 44461; Calidad = "Primera"; Volumen = 40000; Min = 110; Max = 120; Prom = 114 },
    @{ Row = 441; Fecha = 44461; Calidad = "Segunda"; Volumen = 36000; Min = 90;  Max = 95;  Prom = 92 },
    @{ Row = 442; Fecha = 44461; Calidad = "Tercera"; Volumen = 9000;  Min = 65;  Max = 65;  Prom = 65 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $constA
    $ws.Cells.Item($row, 2).Value = $constB
    $ws.Cells.Item($row, 3).Value = $constC
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 5).Value = $constE
    $ws.Cells.Item($row, 6).Value = $constF
    $ws.Cells.Item($row, 7).Value = $constG
    $ws.Cells.Item($row, 8).Value = $constH
    $ws.Cells.Item($row, 9).Value = $r.Calidad
    $ws.Cells.Item($row, 10).Value = $r.Volumen
    $ws.Cells.Item($row, 11).Value = $r.Min
    $ws.Cells.Item($row, 12).Value = $r.Max
    $ws.Cells.Item($row, 13).Value = $r.Prom
    $ws.Cells.Item($row, 14).Value = $constN
    $ws.Cells.Item($row, 15).Value = $constO
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $constQ
    $ws.Cells.Item($row, 18).Value = $constR
}
